$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.127.68'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = '1.989.70'
$ws.Range('E3').Value = '  -2.59%  '
$cell = $ws.Range('D4')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range('E4').Value = '  -0.05%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '244.38'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  -5.22%  '
$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.602'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  -3.66%  '
$ws.Range('E7').Value = '  +0.04%  '
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '54.66'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  -5.18%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '59.52'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  +3.89%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.372'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  -4.11%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0753'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  -5.76%  '
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0985'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  -4.46%  '
$ws.Range('D13').Value = '2.280.57'
$ws.Range('E13').Value = '  -2.59%  '
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '14.01'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  -5.30%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '21.03'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -1.45%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.757'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  -8.02%  '
$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.06'
$cell.Style = $origStyle
$ws.Range('E17').Value = '  -6.16%  '
$ws.Range('D18').Value = '1.996.41'
$ws.Range('E18').Value = '  -2.21%  '
$ws.Range('D19').Value = '37.050.34'
$ws.Range('E19').Value = '  -1.19%  '
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '68.26'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('D21').Value = '0.0₃0810'
$ws.Range('E21').Value = '  -5.53%  '
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '228.94'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  -0.20%  '
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.97'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  -5.22%  '
$ws.Range('E24').Value = '  +0.05%  '
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.43'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  -9.32%  '
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.35'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  -0.10%  '
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '161.13'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  -1.50%  '
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '8.68'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  -5.30%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '19.08'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  -4.61%  '
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.123'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  -10.65%  '
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('E32').Value = '  -3.08%  '
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.42'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  -7.25%  '
$ws.Range('E34').Value = '  -7.57%  '
$ws.Range('E35').Value = '  -6.48%  '
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.35'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  -6.51%  '
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  -0.05%  '
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.78'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  -1.70%  '
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.33'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  -3.85%  '
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.23'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  -2.67%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').Value = '1.423.36'
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('E43').Value = '  -5.25%  '
$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0204'
$cell.Style = $origStyle
$ws.Range('E44').Value = '  -6.56%  '
$ws.Range('E45').Value = '  -8.11%  '
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '88.17'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  -3.72%  '
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '15.35'
$cell.Style = $origStyle
$ws.Range('E47').Value = '  -5.35%  '
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  -5.13%  '
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.87'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  -0.42%  '
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.64'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  -11.05%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.172.45'
$ws.Range('E51').Value = '  -2.61%  '
